$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2364713333333333
$ws.Range("H2").Value = 0.709414
$ws.Range("I2").Value = 0.002249544876489787
$ws.Range("J2").Value = 0.002249544876489787
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 1.770607957836
$ws.Range("R2").Value = 15.935471620524
$ws.Range("S2").Value = 0.0003114259082482917
$ws.Range("T2").Value = 0.0003114259082482917
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2364713333333333
$ws.Range("H3").Value = 0.709414
$ws.Range("I3").Value = 0.002249544876489787
$ws.Range("J3").Value = 0.002249544876489787
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 7.566947405064
$ws.Range("R3").Value = 68.10252664557599
$ws.Range("S3").Value = 0.001330923346334232
$ws.Range("T3").Value = 0.001330923346334232
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2364713333333333
$ws.Range("H4").Value = 0.709414
$ws.Range("I4").Value = 0.002249544876489787
$ws.Range("J4").Value = 0.002249544876489787
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 3.452202824612222
$ws.Range("R4").Value = 31.06982542151
$ws.Range("S4").Value = 0.0006071956219072634
$ws.Range("T4").Value = 0.0006071956219072634
$ws.Range("I5").Value = 0.9862688099613843
$ws.Range("J5").Value = 0.9862688099613843
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 776.2883157983518
$ws.Range("R5").Value = 6986.594842185167
$ws.Range("S5").Value = 0.1365385785939356
$ws.Range("T5").Value = 0.1365385785939356
$ws.Range("I6").Value = 0.9862688099613843
$ws.Range("J6").Value = 0.9862688099613843
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("S6").Value = 0.5835172254874755
$ws.Range("T6").Value = 0.5835172254874755
$ws.Range("I7").Value = 0.9862688099613843
$ws.Range("J7").Value = 0.9862688099613843
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 1513.550588458813
$ws.Range("R7").Value = 13621.95529612932
$ws.Range("S7").Value = 0.2662130058799732
$ws.Range("T7").Value = 0.2662130058799732
$ws.Range("G8").Value = 1.206946333333333
$ws.Range("H8").Value = 3.620839
$ws.Range("I8").Value = 0.01148164516212593
$ws.Range("J8").Value = 0.01148164516212593
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 9.037157918286001
$ws.Range("R8").Value = 81.33442126457399
$ws.Range("S8").Value = 0.001589513421212206
$ws.Range("T8").Value = 0.001589513421212207
$ws.Range("G9").Value = 1.206946333333333
$ws.Range("H9").Value = 3.620839
$ws.Range("I9").Value = 0.01148164516212593
$ws.Range("J9").Value = 0.01148164516212593
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 38.621592293364
$ws.Range("R9").Value = 347.594330640276
$ws.Range("S9").Value = 0.006793013893745395
$ws.Range("T9").Value = 0.006793013893745396
$ws.Range("G10").Value = 1.206946333333333
$ws.Range("H10").Value = 3.620839
$ws.Range("I10").Value = 0.01148164516212593
$ws.Range("J10").Value = 0.01148164516212593
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 17.61999428157056
$ws.Range("R10").Value = 158.579948534135
$ws.Range("S10").Value = 0.00309911784716833
$ws.Range("T10").Value = 0.00309911784716833
